# Rename the embedded logo pictures in the headers/footers:
#   - The two Pearson logo pictures (in the "first page" footer and the
#     "default" footer) go from "image2.png" to "image1.png".
#   - The BTec logo picture (in the "first page" header) goes from
#     "image1.jpg" to "image2.jpg".
#
# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer "first page" (footer1.xml) - Pearson logo, id=3
$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Exists -and $ftrFirst.Range.InlineShapes.Count -ge 1) {
    $pearson1 = $ftrFirst.Range.InlineShapes.Item(1)
    $pearson1.Name = "image1.png"
}

# Footer "default/primary" (footer2.xml) - Pearson logo, id=2
$ftrPrimary = $sec.Footers.Item(1)
if ($ftrPrimary.Exists -and $ftrPrimary.Range.InlineShapes.Count -ge 1) {
    $pearson2 = $ftrPrimary.Range.InlineShapes.Item(1)
    $pearson2.Name = "image1.png"
}

# Header "first page" (header1.xml) - BTec logo, id=1
$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Exists -and $hdrFirst.Range.InlineShapes.Count -ge 1) {
    $btec = $hdrFirst.Range.InlineShapes.Item(1)
    $btec.Name = "image2.jpg"
}
